$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '61.494.53'
$ws.Range('E2').Value = '  +1.18%  '

$ws.Range('D3').Value = '3.386.55'
$ws.Range('E3').Value = '  +0.89%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').Value = '575.96'
$ws.Range('E5').Value = '  +1.00%  '

$ws.Range('D6').Value = '136.89'
$ws.Range('E6').Value = '  +2.31%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').Value = '3.386.64'
$ws.Range('E8').Value = '  +0.94%  '

$ws.Range('E9').Value = '  -0.63%  '

$ws.Range('D10').Value = '7.50'
$ws.Range('E10').Value = '  -1.24%  '

$ws.Range('D11').Value = '0.125'
$ws.Range('E11').Value = '  +2.96%  '

$ws.Range('D12').Value = '0.389'
$ws.Range('E12').Value = '  +0.70%  '

$ws.Range('D13').Value = '3.966.29'
$ws.Range('E13').Value = '  +0.84%  '

$ws.Range('E14').Value = '  +2.83%  '

$ws.Range('D15').Value = '0.0000176'
$ws.Range('E15').Value = '  +2.81%  '

$ws.Range('D16').Value = '3.384.48'
$ws.Range('E16').Value = '  +0.74%  '

$ws.Range('D17').Value = '25.79'
$ws.Range('E17').Value = '  +3.33%  '

$ws.Range('D18').Value = '61.591.05'
$ws.Range('E18').Value = '  +1.14%  '

$ws.Range('D19').Value = '14.18'
$ws.Range('E19').Value = '  +2.37%  '

$ws.Range('D20').Value = '5.86'
$ws.Range('E20').Value = '  +2.34%  '

$ws.Range('D21').Value = '9.43'
$ws.Range('E21').Value = '  +0.61%  '

$ws.Range('D22').Value = '376.17'
$ws.Range('E22').Value = '  +1.12%  '

$ws.Range('D23').Value = '0.557'
$ws.Range('E23').Value = '  -2.45%  '

$ws.Range('D24').Value = '3.531.93'
$ws.Range('E24').Value = '  +1.12%  '

$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.03%  '

$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '0.0000127'
$ws.Range('E26').Value = '  +9.01%  '

$ws.Range('D27').Value = '71.24'
$ws.Range('E27').Value = '  +1.04%  '

$ws.Range('D28').Value = '1.69'
$ws.Range('E28').Value = '  +1.80%  '

$ws.Range('D29').Value = '7.51'
$ws.Range('E29').Value = '  -1.48%  '

$ws.Range('E30').Value = '  -0.11%  '

$ws.Range('E31').Value = '  +4.89%  '

$ws.Range('D32').Value = '8.23'
$ws.Range('E32').Value = '  +2.05%  '

$ws.Range('D33').Value = '2.17'
$ws.Range('E33').Value = '  +2.01%  '

$ws.Range('E34').Value = '  +0.05%  '

$ws.Range('E35').Value = '  +0.78%  '

$ws.Range('D36').Value = '5.32'
$ws.Range('E36').Value = '  -3.65%  '

$ws.Range('E37').Value = '  +0.23%  '

$ws.Range('D38').Value = '6.83'
$ws.Range('E38').Value = '  -0.81%  '

$ws.Range('D39').Value = '165.38'
$ws.Range('E39').Value = '  +1.97%  '

$ws.Range('D40').Value = '0.0780'
$ws.Range('E40').Value = '  +0.08%  '

$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D41').Value = '0.781'
$ws.Range('E41').Value = '  +3.35%  '

$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.10%  '

$ws.Range('B43').Value = 'ONDO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D43').Value = '1.23'
$ws.Range('E43').Value = '  +2.81%  '

$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').Value = '1.72'
$ws.Range('E44').Value = '  +8.20%  '

$ws.Range('D45').Value = '25.02'
$ws.Range('E45').Value = '  +9.60%  '

$ws.Range('D46').Value = '4.41'
$ws.Range('E46').Value = '  +0.94%  '

$ws.Range('D47').Value = '41.44'
$ws.Range('E47').Value = '  +0.53%  '

$ws.Range('D48').Value = '6.85'
$ws.Range('E48').Value = '  -1.22%  '

$ws.Range('D49').Value = '22.58'
$ws.Range('E49').Value = '  -2.82%  '

$ws.Range('D50').Value = '2.333.82'
$ws.Range('E50').Value = '  +5.31%  '

$ws.Range('B51').Value = 'LidoDAOToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D51').Value = '2.38'
$ws.Range('E51').Value = '  -0.84%  '

$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D12').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D17').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').Style = "Normal"
$ws.Range('D29').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').Style = "Normal"
$ws.Range('D51').Style = "Normal"

